$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -5.743499999999998
$ws.Range("C9").Value = -11.9645
$ws.Range("C18").Value = -14.4533
$ws.Range("C20").Value = -13.71099999999998
